$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(82, 8).Value = 1678.3  # ALC!H82 (was 2429.375)
$ws.Cells.Item(82, 9).Value = 683.2857  # ALC!I82 (was 858.75)
$ws.Cells.Item(82, 11).Value = 2049.8571  # ALC!K82 (was 2576.25)
$ws.Cells.Item(82, 13).Value = -1643.8571  # ALC!M82 (was -2170.25)
$ws.Cells.Item(85, 8).Value = 1678.3  # ALC!H85 (was 2429.375)
$ws.Cells.Item(85, 9).Value = 683.2857  # ALC!I85 (was 858.75)
$ws.Cells.Item(85, 11).Value = 2049.8571  # ALC!K85 (was 2576.25)
$ws.Cells.Item(85, 13).Value = -645.8571000000002  # ALC!M85 (was -1172.25)
$ws.Cells.Item(86, 8).Value = 81088.60000000001  # ALC!H86 (was 42671.418)
$ws.Cells.Item(86, 9).Value = 1574.3334  # ALC!I86 (was 13613.375)
$ws.Cells.Item(86, 10).Value = 200360  # ALC!J86 (was 100787.5)
$ws.Cells.Item(86, 11).Value = 1574.3334  # ALC!K86 (was 13613.375)
$ws.Cells.Item(86, 12).Value = 200360  # ALC!L86 (was 100787.5)
$ws.Cells.Item(86, 13).Value = -451.3334  # ALC!M86 (was -12490.375)
$ws.Cells.Item(86, 14).Value = -202606  # ALC!N86 (was -103033.5)
$ws.Cells.Item(89, 8).Value = 81088.60000000001  # ALC!H89 (was 42671.418)
$ws.Cells.Item(89, 9).Value = 1574.3334  # ALC!I89 (was 13613.375)
$ws.Cells.Item(89, 10).Value = 200360  # ALC!J89 (was 100787.5)
$ws.Cells.Item(89, 11).Value = 7871.666999999999  # ALC!K89 (was 68066.875)
$ws.Cells.Item(89, 12).Value = 1001800  # ALC!L89 (was 503937.5)
$ws.Cells.Item(89, 13).Value = -2255.666999999999  # ALC!M89 (was -62450.875)
$ws.Cells.Item(89, 14).Value = -1013032  # ALC!N89 (was -515169.5)
$ws.Cells.Item(98, 8).Value = 764.0454999999999  # ALC!H98 (was 790.15)
$ws.Cells.Item(98, 9).Value = 769.4737  # ALC!I98 (was 806.41174)
$ws.Cells.Item(98, 10).Value = 729.6667  # ALC!J98 (was 698)
$ws.Cells.Item(98, 11).Value = 769.4737  # ALC!K98 (was 806.41174)
$ws.Cells.Item(98, 12).Value = 729.6667  # ALC!L98 (was 698)
$ws.Cells.Item(98, 13).Value = 728.5263  # ALC!M98 (was 691.58826)
$ws.Cells.Item(98, 14).Value = -3725.6667  # ALC!N98 (was -3694)
$ws.Cells.Item(122, 8).Value = 764.0454999999999  # ALC!H122 (was 790.15)
$ws.Cells.Item(122, 9).Value = 769.4737  # ALC!I122 (was 806.41174)
$ws.Cells.Item(122, 10).Value = 729.6667  # ALC!J122 (was 698)
$ws.Cells.Item(122, 11).Value = 2308.4211  # ALC!K122 (was 2419.23522)
$ws.Cells.Item(122, 12).Value = 2189.0001  # ALC!L122 (was 2094)
$ws.Cells.Item(122, 13).Value = 141.5789  # ALC!M122 (was 30.76477999999997)
$ws.Cells.Item(122, 14).Value = -7089.0001  # ALC!N122 (was -6994)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 840.1070999999999  # ARM!H2 (was 901.6957)
$ws.Cells.Item(2, 9).Value = 802.087  # ARM!I2 (was 850.7368)
$ws.Cells.Item(2, 10).Value = 1015  # ARM!J2 (was 1143.75)
$ws.Cells.Item(2, 11).Value = 802.087  # ARM!K2 (was 850.7368)
$ws.Cells.Item(2, 12).Value = 1015  # ARM!L2 (was 1143.75)
$ws.Cells.Item(2, 13).Value = -689.087  # ARM!M2 (was -737.7368)
$ws.Cells.Item(2, 14).Value = -1241  # ARM!N2 (was -1369.75)
$ws.Cells.Item(32, 8).Value = 5109.9663  # ARM!H32 (was 5757.4614)
$ws.Cells.Item(32, 9).Value = 2841.2173  # ARM!I32 (was 3258.6924)
$ws.Cells.Item(32, 10).Value = 13137.846  # ARM!J32 (was 13253.77)
$ws.Cells.Item(32, 11).Value = 2841.2173  # ARM!K32 (was 3258.6924)
$ws.Cells.Item(32, 12).Value = 13137.846  # ARM!L32 (was 13253.77)
$ws.Cells.Item(32, 13).Value = -2554.2173  # ARM!M32 (was -2971.6924)
$ws.Cells.Item(32, 14).Value = -13711.846  # ARM!N32 (was -13827.77)
$ws.Cells.Item(35, 8).Value = 3037  # ARM!H35 (was 0)
$ws.Cells.Item(35, 9).Value = 3037  # ARM!I35 (was 0)
$ws.Cells.Item(35, 11).Value = 3037  # ARM!K35 (was 0)
$ws.Cells.Item(35, 13).Value = -2631  # ARM!M35 (was None)
$ws.Cells.Item(36, 8).Value = 6750  # ARM!H36 (was 8005.2)
$ws.Cells.Item(36, 9).Value = 7333.3335  # ARM!I36 (was 8005.2)
$ws.Cells.Item(36, 10).Value = 5000  # ARM!J36 (was 0)
$ws.Cells.Item(36, 11).Value = 7333.3335  # ARM!K36 (was 8005.2)
$ws.Cells.Item(36, 12).Value = 5000  # ARM!L36 (was 0)
$ws.Cells.Item(36, 13).Value = -6987.3335  # ARM!M36 (was -7659.2)
$ws.Cells.Item(36, 14).Value = -5692  # ARM!N36 (was None)
$ws.Cells.Item(51, 8).Value = 50000  # ARM!H51 (was 76547)
$ws.Cells.Item(51, 10).Value = 50000  # ARM!J51 (was 76547)
$ws.Cells.Item(51, 12).Value = 50000  # ARM!L51 (was 76547)
$ws.Cells.Item(51, 14).Value = -51512  # ARM!N51 (was -78059)
$ws.Cells.Item(93, 8).Value = 63448  # ARM!H93 (was 0)
$ws.Cells.Item(93, 10).Value = 63448  # ARM!J93 (was 0)
$ws.Cells.Item(93, 12).Value = 63448  # ARM!L93 (was 0)
$ws.Cells.Item(93, 14).Value = -68440  # ARM!N93 (was None)
$ws.Cells.Item(98, 8).Value = 9000  # ARM!H98 (was 0)
$ws.Cells.Item(98, 10).Value = 9000  # ARM!J98 (was 0)
$ws.Cells.Item(98, 12).Value = 9000  # ARM!L98 (was 0)
$ws.Cells.Item(98, 14).Value = -14990  # ARM!N98 (was None)
$ws.Cells.Item(116, 8).Value = 840.1070999999999  # ARM!H116 (was 901.6957)
$ws.Cells.Item(116, 9).Value = 802.087  # ARM!I116 (was 850.7368)
$ws.Cells.Item(116, 10).Value = 1015  # ARM!J116 (was 1143.75)
$ws.Cells.Item(116, 11).Value = 802.087  # ARM!K116 (was 850.7368)
$ws.Cells.Item(116, 12).Value = 1015  # ARM!L116 (was 1143.75)
$ws.Cells.Item(116, 13).Value = 1491.913  # ARM!M116 (was 1443.2632)
$ws.Cells.Item(116, 14).Value = -5603  # ARM!N116 (was -5731.75)
$ws.Cells.Item(122, 8).Value = 1515.325  # ARM!H122 (was 1700.069)
$ws.Cells.Item(122, 9).Value = 1723.7059  # ARM!I122 (was 2026.5454)
$ws.Cells.Item(122, 10).Value = 1361.3043  # ARM!J122 (was 1500.5555)
$ws.Cells.Item(122, 11).Value = 5171.1177  # ARM!K122 (was 6079.6362)
$ws.Cells.Item(122, 12).Value = 4083.9129  # ARM!L122 (was 4501.666499999999)
$ws.Cells.Item(122, 13).Value = -2721.1177  # ARM!M122 (was -3629.6362)
$ws.Cells.Item(122, 14).Value = -8983.912899999999  # ARM!N122 (was -9401.666499999999)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 840.1070999999999  # BSM!H3 (was 901.6957)
$ws.Cells.Item(3, 9).Value = 802.087  # BSM!I3 (was 850.7368)
$ws.Cells.Item(3, 10).Value = 1015  # BSM!J3 (was 1143.75)
$ws.Cells.Item(3, 11).Value = 802.087  # BSM!K3 (was 850.7368)
$ws.Cells.Item(3, 12).Value = 1015  # BSM!L3 (was 1143.75)
$ws.Cells.Item(3, 13).Value = -688.087  # BSM!M3 (was -736.7368)
$ws.Cells.Item(3, 14).Value = -1243  # BSM!N3 (was -1371.75)
$ws.Cells.Item(36, 8).Value = 1000  # BSM!H36 (was 2999)
$ws.Cells.Item(36, 9).Value = 1000  # BSM!I36 (was 2999)
$ws.Cells.Item(36, 11).Value = 1000  # BSM!K36 (was 2999)
$ws.Cells.Item(36, 13).Value = -466  # BSM!M36 (was -2465)
$ws.Cells.Item(94, 8).Value = 629  # BSM!H94 (was 631.6667)
$ws.Cells.Item(94, 9).Value = 698.3333  # BSM!I94 (was 685)
$ws.Cells.Item(94, 11).Value = 698.3333  # BSM!K94 (was 685)
$ws.Cells.Item(94, 13).Value = -247.3333  # BSM!M94 (was -234)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 19613166  # CRP!H31 (was 21282314)
$ws.Cells.Item(31, 9).Value = 38462396  # CRP!I31 (was 43479156)
$ws.Cells.Item(31, 10).Value = 9966.120000000001  # CRP!J31 (was 10339.708)
$ws.Cells.Item(31, 11).Value = 38462396  # CRP!K31 (was 43479156)
$ws.Cells.Item(31, 12).Value = 9966.120000000001  # CRP!L31 (was 10339.708)
$ws.Cells.Item(31, 13).Value = -38462101  # CRP!M31 (was -43478861)
$ws.Cells.Item(31, 14).Value = -10556.12  # CRP!N31 (was -10929.708)
$ws.Cells.Item(34, 8).Value = 19613166  # CRP!H34 (was 21282314)
$ws.Cells.Item(34, 9).Value = 38462396  # CRP!I34 (was 43479156)
$ws.Cells.Item(34, 10).Value = 9966.120000000001  # CRP!J34 (was 10339.708)
$ws.Cells.Item(34, 11).Value = 38462396  # CRP!K34 (was 43479156)
$ws.Cells.Item(34, 12).Value = 9966.120000000001  # CRP!L34 (was 10339.708)
$ws.Cells.Item(34, 13).Value = -38462194  # CRP!M34 (was -43478954)
$ws.Cells.Item(34, 14).Value = -10370.12  # CRP!N34 (was -10743.708)
$ws.Cells.Item(109, 8).Value = 37293.332  # CRP!H109 (was 37926.668)
$ws.Cells.Item(109, 10).Value = 37293.332  # CRP!J109 (was 37926.668)
$ws.Cells.Item(109, 12).Value = 37293.332  # CRP!L109 (was 37926.668)
$ws.Cells.Item(109, 14).Value = -39373.332  # CRP!N109 (was -40006.668)
$ws.Cells.Item(122, 8).Value = 2000  # CRP!H122 (was 1247.6364)
$ws.Cells.Item(122, 9).Value = 2000  # CRP!I122 (was 1072.4)
$ws.Cells.Item(122, 10).Value = 0  # CRP!J122 (was 3000)
$ws.Cells.Item(122, 11).Value = 6000  # CRP!K122 (was 3217.2)
$ws.Cells.Item(122, 12).Value = 0  # CRP!L122 (was 9000)
$ws.Cells.Item(122, 13).Value = -3550  # CRP!M122 (was -767.2000000000003)
$ws.Cells.Item(122, 14).ClearContents()  # CRP!N122 (was -13900)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 4782.75  # CUL!H62 (was 4246.2)
$ws.Cells.Item(62, 10).Value = 6123  # CUL!J62 (was 5117.25)
$ws.Cells.Item(62, 12).Value = 18369  # CUL!L62 (was 15351.75)
$ws.Cells.Item(62, 14).Value = -19741  # CUL!N62 (was -16723.75)
$ws.Cells.Item(65, 8).Value = 4782.75  # CUL!H65 (was 4246.2)
$ws.Cells.Item(65, 10).Value = 6123  # CUL!J65 (was 5117.25)
$ws.Cells.Item(65, 12).Value = 55107  # CUL!L65 (was 46055.25)
$ws.Cells.Item(65, 14).Value = -61971  # CUL!N65 (was -52919.25)
$ws.Cells.Item(107, 8).Value = 807.1  # CUL!H107 (was 760.41174)
$ws.Cells.Item(107, 9).Value = 849.2727  # CUL!I107 (was 884.8095)
$ws.Cells.Item(107, 10).Value = 755.55554  # CUL!J107 (was 673.3333)
$ws.Cells.Item(107, 11).Value = 2547.8181  # CUL!K107 (was 2654.4285)
$ws.Cells.Item(107, 12).Value = 2266.66662  # CUL!L107 (was 2019.9999)
$ws.Cells.Item(107, 13).Value = -627.8181  # CUL!M107 (was -734.4285)
$ws.Cells.Item(107, 14).Value = -6106.66662  # CUL!N107 (was -5859.9999)
$ws.Cells.Item(131, 8).Value = 992.5876500000001  # CUL!H131 (was 992.0106)
$ws.Cells.Item(131, 10).Value = 992.5876500000001  # CUL!J131 (was 992.0106)
$ws.Cells.Item(131, 12).Value = 2977.76295  # CUL!L131 (was 2976.0318)
$ws.Cells.Item(131, 14).Value = -13057.76295  # CUL!N131 (was -13056.0318)
$ws.Cells.Item(136, 8).Value = 1879.375  # CUL!H136 (was 1515.909)
$ws.Cells.Item(136, 9).Value = 1785.5555  # CUL!I136 (was 1515.909)
$ws.Cells.Item(136, 10).Value = 2000  # CUL!J136 (was 0)
$ws.Cells.Item(136, 11).Value = 5356.666499999999  # CUL!K136 (was 4547.727000000001)
$ws.Cells.Item(136, 12).Value = 6000  # CUL!L136 (was 0)
$ws.Cells.Item(136, 13).Value = -256.6664999999994  # CUL!M136 (was 552.2729999999992)
$ws.Cells.Item(136, 14).Value = -16200  # CUL!N136 (was None)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10460.0625  # GSM!H70 (was 10097.529)
$ws.Cells.Item(70, 10).Value = 4579.2856  # GSM!J70 (was 4544)
$ws.Cells.Item(70, 12).Value = 4579.2856  # GSM!L70 (was 4544)
$ws.Cells.Item(70, 14).Value = -5119.2856  # GSM!N70 (was -5084)
$ws.Cells.Item(73, 8).Value = 10460.0625  # GSM!H73 (was 10097.529)
$ws.Cells.Item(73, 10).Value = 4579.2856  # GSM!J73 (was 4544)
$ws.Cells.Item(73, 12).Value = 4579.2856  # GSM!L73 (was 4544)
$ws.Cells.Item(73, 14).Value = -6451.2856  # GSM!N73 (was -6416)
$ws.Cells.Item(99, 8).Value = 9910.4  # GSM!H99 (was 9300.454)
$ws.Cells.Item(99, 9).Value = 2847.5  # GSM!I99 (was 2898)
$ws.Cells.Item(99, 11).Value = 2847.5  # GSM!K99 (was 2898)
$ws.Cells.Item(99, 13).Value = -601.5  # GSM!M99 (was -652)
$ws.Cells.Item(122, 8).Value = 2134.1086  # GSM!H122 (was 2063.578)
$ws.Cells.Item(122, 9).Value = 1807.1316  # GSM!I122 (was 1905.5161)
$ws.Cells.Item(122, 10).Value = 3687.25  # GSM!J122 (was 2413.5715)
$ws.Cells.Item(122, 11).Value = 5421.3948  # GSM!K122 (was 5716.5483)
$ws.Cells.Item(122, 12).Value = 11061.75  # GSM!L122 (was 7240.7145)
$ws.Cells.Item(122, 13).Value = -2971.3948  # GSM!M122 (was -3266.5483)
$ws.Cells.Item(122, 14).Value = -15961.75  # GSM!N122 (was -12140.7145)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 49498.91  # LTW!H40 (was 49543.91)
$ws.Cells.Item(40, 9).Value = 59409.777  # LTW!I40 (was 59464.777)
$ws.Cells.Item(40, 11).Value = 59409.777  # LTW!K40 (was 59464.777)
$ws.Cells.Item(40, 13).Value = -59273.777  # LTW!M40 (was -59328.777)
